$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update Bali's (row 8) raw scores per the diff
$ws.Range("B8").Value = 95
$ws.Range("E8").Value = 97
$ws.Range("F8").Value = 96

# Reflect the updated selection/scroll position recorded in the sheet view
$ws.Range("L22").Select() | Out-Null
